$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 updates
$ws.Range("G11").Value = 1.91
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.95
$ws.Range("L11").Value = 1.39
$ws.Range("M11").Value = 2.55
$ws.Range("N11").Value = 2.15
$ws.Range("O11").Value = 1.55
$ws.Range("P11").Value = 1.45
$ws.Range("Q11").Value = 2.4
$ws.Range("R11").Value = 1.93
$ws.Range("S11").Value = 1.7
$ws.Range("T11").Value = 5.8
$ws.Range("U11").Value = 8
$ws.Range("V11").Value = 8.75
$ws.Range("W11").Value = 16
$ws.Range("X11").Value = 17.5
$ws.Range("Y11").Value = 35
$ws.Range("Z11").Value = 7.6
$ws.Range("AB11").Value = 17.5
$ws.Range("AC11").Value = 100
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 9.5
$ws.Range("AF11").Value = 20
$ws.Range("AG11").Value = 13.5
$ws.Range("AH11").Value = 60
$ws.Range("AJ11").Value = 55

# Row 13 updates
$ws.Range("L13").Value = 1.33
$ws.Range("M13").Value = 3.25

# Row 14 updates
$ws.Range("G14").Value = 1.6
$ws.Range("H14").Value = 3.85
$ws.Range("I14").Value = 4.7
$ws.Range("K14").Value = 8.25
$ws.Range("L14").Value = 1.23
$ws.Range("M14").Value = 3.75
$ws.Range("N14").Value = 1.7
$ws.Range("O14").Value = 2.05
$ws.Range("P14").Value = 1.34
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 1.75
$ws.Range("S14").Value = 1.98
$ws.Range("T14").Value = 7.8
$ws.Range("U14").Value = 8
$ws.Range("V14").Value = 8
$ws.Range("W14").Value = 12
$ws.Range("X14").Value = 12
$ws.Range("Y14").Value = 23
$ws.Range("Z14").Value = 8.25
$ws.Range("AA14").Value = 7.8
$ws.Range("AB14").Value = 15.5
$ws.Range("AC14").Value = 65
$ws.Range("AD14").Value = 450
$ws.Range("AE14").Value = 14.5
$ws.Range("AF14").Value = 29
$ws.Range("AG14").Value = 15.5
$ws.Range("AH14").Value = 80
